$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C, rows 2 through 18 hold a "Förändrad" (changed/updated) date.
# The diff bumps the serial date value from 45189 to 45190 (2023-09-20 -> 2023-09-21)
# for every row in that range, leaving everything else untouched.
for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value2 = 45190
    }
}
